# Auto-generated edit script applying the Marilith_Profits market-data refresh diff
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 102794.8
$ws.Range("I74").Value = 3493.5
$ws.Range("K74").Value = 3493.5
$ws.Range("M74").Value = -2557.5
$ws.Range("H77").Value = 102794.8
$ws.Range("I77").Value = 3493.5
$ws.Range("K77").Value = 17467.5
$ws.Range("M77").Value = -12787.5
$ws.Range("H80").Value = 634.8421
$ws.Range("J80").Value = 624.4545000000001
$ws.Range("L80").Value = 1873.3635
$ws.Range("N80").Value = -3869.3635
$ws.Range("H83").Value = 634.8421
$ws.Range("J83").Value = 624.4545000000001
$ws.Range("L83").Value = 5620.0905
$ws.Range("N83").Value = -15604.0905
$ws.Range("H88").Value = 2550
$ws.Range("I88").Value = 3750
$ws.Range("J88").Value = 750
$ws.Range("K88").Value = 3750
$ws.Range("L88").Value = 750
$ws.Range("M88").Value = -3344
$ws.Range("N88").Value = -1562
$ws.Range("H91").Value = 2550
$ws.Range("I91").Value = 3750
$ws.Range("J91").Value = 750
$ws.Range("K91").Value = 3750
$ws.Range("L91").Value = 750
$ws.Range("M91").Value = -2346
$ws.Range("N91").Value = -3558
$ws.Range("H111").Value = 751.6667
$ws.Range("I111").Value = 682.2857
$ws.Range("K111").Value = 2046.8571
$ws.Range("M111").Value = 1020.1429
$ws.Range("H112").Value = 1953.4286
$ws.Range("I112").Value = 1449.6666
$ws.Range("J112").Value = 2090.818
$ws.Range("K112").Value = 4348.9998
$ws.Range("L112").Value = 6272.454000000001
$ws.Range("M112").Value = -3240.9998
$ws.Range("N112").Value = -8488.454000000002
$ws.Range("H127").Value = 1482
$ws.Range("I127").Value = 1178.4
$ws.Range("K127").Value = 3535.2
$ws.Range("M127").Value = 1424.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 470.42856
$ws.Range("I2").Value = 470.5
$ws.Range("J2").Value = 470
$ws.Range("K2").Value = 470.5
$ws.Range("L2").Value = 470
$ws.Range("M2").Value = -357.5
$ws.Range("N2").Value = -696
$ws.Range("H32").Value = 1955.1034
$ws.Range("I32").Value = 1667.7858
$ws.Range("K32").Value = 1667.7858
$ws.Range("M32").Value = -1380.7858
$ws.Range("H97").Value = 787.3333
$ws.Range("J97").Value = 1003.6667
$ws.Range("L97").Value = 1003.6667
$ws.Range("N97").Value = -1995.6667
$ws.Range("H116").Value = 470.42856
$ws.Range("I116").Value = 470.5
$ws.Range("J116").Value = 470
$ws.Range("K116").Value = 470.5
$ws.Range("L116").Value = 470
$ws.Range("M116").Value = 1823.5
$ws.Range("N116").Value = -5058

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 470.42856
$ws.Range("I3").Value = 470.5
$ws.Range("J3").Value = 470
$ws.Range("K3").Value = 470.5
$ws.Range("L3").Value = 470
$ws.Range("M3").Value = -356.5
$ws.Range("N3").Value = -698
$ws.Range("H11").Value = 226
$ws.Range("I11").Value = 226
$ws.Range("K11").Value = 226
$ws.Range("M11").Value = -86
$ws.Range("H20").Value = 1284.5
$ws.Range("I20").Value = 985
$ws.Range("K20").Value = 985
$ws.Range("M20").Value = -738
$ws.Range("H94").Value = 1009.38464
$ws.Range("I94").Value = 960.1667
$ws.Range("K94").Value = 960.1667
$ws.Range("M94").Value = -509.1667
$ws.Range("H105").Value = 4605.9
$ws.Range("I105").Value = 4605.9
$ws.Range("K105").Value = 4605.9
$ws.Range("M105").Value = -2858.9
$ws.Range("H134").Value = 8536.388999999999
$ws.Range("I134").Value = 8536.388999999999
$ws.Range("K134").Value = 25609.167
$ws.Range("M134").Value = -23074.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 702.3333
$ws.Range("I22").Value = 364.2
$ws.Range("K22").Value = 364.2
$ws.Range("M22").Value = -14.19999999999999
$ws.Range("H99").Value = 4430.1113
$ws.Range("I99").Value = 3616.3333
$ws.Range("J99").Value = 4837
$ws.Range("K99").Value = 3616.3333
$ws.Range("L99").Value = 4837
$ws.Range("M99").Value = -2118.3333
$ws.Range("N99").Value = -7833
$ws.Range("H107").Value = 690.4706
$ws.Range("I107").Value = 645.3333
$ws.Range("J107").Value = 798.8
$ws.Range("K107").Value = 645.3333
$ws.Range("L107").Value = 798.8
$ws.Range("M107").Value = 1274.6667
$ws.Range("N107").Value = -4638.8
$ws.Range("H126").Value = 4430.1113
$ws.Range("I126").Value = 3616.3333
$ws.Range("J126").Value = 4837
$ws.Range("K126").Value = 10848.9999
$ws.Range("L126").Value = 14511
$ws.Range("M126").Value = -8378.999899999999
$ws.Range("N126").Value = -19451

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99900
$ws.Range("J37").Value = 99900
$ws.Range("L37").Value = 299700
$ws.Range("N37").Value = -299924
$ws.Range("H129").Value = 731.5
$ws.Range("J129").Value = 913
$ws.Range("L129").Value = 2739
$ws.Range("N129").Value = -12739
$ws.Range("H132").Value = 5454
$ws.Range("I132").Value = 5454
$ws.Range("K132").Value = 49086
$ws.Range("M132").Value = -46556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1507.375
$ws.Range("I102").Value = 1508.4286
$ws.Range("K102").Value = 1508.4286
$ws.Range("M102").Value = 113.5714
$ws.Range("H107").Value = 171.63637
$ws.Range("J107").Value = 90
$ws.Range("L107").Value = 90
$ws.Range("N107").Value = -3930

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1998.8
$ws.Range("I22").Value = 1997.5
$ws.Range("K22").Value = 1997.5
$ws.Range("M22").Value = -1702.5
$ws.Range("H27").Value = 1998.8
$ws.Range("I27").Value = 1997.5
$ws.Range("K27").Value = 1997.5
$ws.Range("M27").Value = -1890.5
$ws.Range("H40").Value = 7946
$ws.Range("I40").Value = 7466.3335
$ws.Range("J40").Value = 8425.666999999999
$ws.Range("K40").Value = 7466.3335
$ws.Range("L40").Value = 8425.666999999999
$ws.Range("M40").Value = -7330.3335
$ws.Range("N40").Value = -8697.666999999999
$ws.Range("H55").Value = 212.09091
$ws.Range("J55").Value = 248
$ws.Range("L55").Value = 248
$ws.Range("N55").Value = -594
$ws.Range("H61").Value = 2335.2856
$ws.Range("J61").Value = 1500
$ws.Range("L61").Value = 1500
$ws.Range("N61").Value = -1904
$ws.Range("H93").Value = 1646.625
$ws.Range("I93").Value = 1567.5714
$ws.Range("K93").Value = 1567.5714
$ws.Range("M93").Value = -319.5714
$ws.Range("H107").Value = 7597.8
$ws.Range("I107").Value = 7597.8
$ws.Range("K107").Value = 7597.8
$ws.Range("M107").Value = -5677.8
$ws.Range("H113").Value = 2335.2856
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840
$ws.Range("H136").Value = 7999.6665
$ws.Range("I136").Value = 7499.5
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 22498.5
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -19948.5
$ws.Range("N136").Value = -32100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1028.5714
$ws.Range("I100").Value = 740
$ws.Range("K100").Value = 1480
$ws.Range("M100").Value = -939
